# Change network line settings:
# The R + L + (C // G) branch has been changed to (R + L)//C//G branch.
#
# Self branches (From bus == To bus) used to model the shunt G/C with R=0, wL=0
# (i.e. R+L short-circuited) and G=0 (i.e. G open). Now the series R+L path is
# opened (set to "inf") so the branch is purely the parallel C (and G) path.
# The former self branches (rows 10/11, real From/To branches) used G="inf"
# (open) which now becomes 0 (closed / not open) since G no longer needs to be
# forced open on those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NetworkLine")

# Mutual branches 1-2 and 3-4: G (pu) column goes from "inf" to 0
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 0

# Self branches (bus 1..4): R (pu) and wL (pu) columns go from 0 to "inf"
$ws.Range("C12").Value = "inf"
$ws.Range("D12").Value = "inf"
$ws.Range("C13").Value = "inf"
$ws.Range("D13").Value = "inf"
$ws.Range("C14").Value = "inf"
$ws.Range("D14").Value = "inf"
$ws.Range("C15").Value = "inf"
$ws.Range("D15").Value = "inf"

# Make NetworkLine the active sheet/tab with F11 selected, matching the
# recorded UI state after the edit (previously the Apparatus sheet was active
# with A2:XFD2 selected).
$null = $ws.Activate()
$null = $ws.Range("F11").Select()
